# Attendance record sheet with a data row for Linda Anderson.
# Mirrors: new "AttendanceRecord" worksheet appended after "PerfomanceReview",
# containing the "employeeName" header and "Linda Anderson" value, becoming
# the active/selected sheet (tab position 3, cell A5 selected), while the
# previously active "PerfomanceReview" sheet's selection moves to F2.

$wb = $excel.ActiveWorkbook

# Add the new worksheet as the very last tab (after the current last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "AttendanceRecord"

# Populate the attendance record data.
$newSheet.Range("A1").Value = "employeeName"
$newSheet.Range("A2").Value = "Linda Anderson"

# Size column A to fit its contents (as Excel would after typing the data).
$newSheet.Columns.Item(1).AutoFit()

# Move the old active sheet's selection off of E2 onto F2 (it is no longer
# the selected/active tab once the new sheet takes over).
$perfSheet = $wb.Worksheets.Item("PerfomanceReview")
$perfSheet.Range("F2").Select()

# Finally, select A5 on the new sheet so it becomes the active tab/cell.
$newSheet.Range("A5").Select()
